# Auto-generated PowerShell Excel COM-interop script
# Applies cell value updates to Sheet1 per the target diff
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Ensure column D keeps its original text formatting so numeric-looking
# strings (e.g. '583.51', '0.0000180') are not coerced into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '64.284.24'
$ws.Range('E2').Value = '  -0.74%  '
$ws.Range('D3').Value = '3.506.92'
$ws.Range('E3').Value = '  +0.00%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '583.51'
$ws.Range('E5').Value = '  -0.50%  '
$ws.Range('D6').Value = '134.79'
$ws.Range('E6').Value = '  +1.39%  '
$ws.Range('D7').Value = '3.507.87'
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  -0.11%  '
$ws.Range('E10').Value = '  +0.07%  '
$ws.Range('D11').Value = '7.10'
$ws.Range('D12').Value = '0.376'
$ws.Range('E12').Value = '  -2.68%  '
$ws.Range('D13').Value = '4.102.26'
$ws.Range('E13').Value = '  -0.05%  '
$ws.Range('D14').Value = '27.14'
$ws.Range('E14').Value = '  -2.23%  '
$ws.Range('D15').Value = '0.0000180'
$ws.Range('E15').Value = '  -0.28%  '
$ws.Range('E16').Value = '  +1.12%  '
$ws.Range('D17').Value = '3.503.25'
$ws.Range('E17').Value = '  -0.13%  '
$ws.Range('D18').Value = '64.300.86'
$ws.Range('E18').Value = '  -0.74%  '
$ws.Range('D19').Value = '9.76'
$ws.Range('E19').Value = '  -2.41%  '
$ws.Range('D20').Value = '13.86'
$ws.Range('E20').Value = '  -2.80%  '
$ws.Range('E21').Value = '  -1.57%  '
$ws.Range('D22').Value = '384.11'
$ws.Range('E22').Value = '  -1.68%  '
$ws.Range('E23').Value = '  -1.48%  '
$ws.Range('D24').Value = '3.646.26'
$ws.Range('E24').Value = '  -0.07%  '
$ws.Range('D25').Value = '73.98'
$ws.Range('E25').Value = '  -0.39%  '
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('D27').Value = '5.74'
$ws.Range('E27').Value = '  +0.04%  '
$ws.Range('D28').Value = '0.0000114'
$ws.Range('E28').Value = '  +4.34%  '
$ws.Range('E29').Value = '  -0.43%  '
$ws.Range('D30').Value = '7.60'
$ws.Range('E30').Value = '  +2.13%  '
$ws.Range('E31').Value = '  +0.03%  '
$ws.Range('D32').Value = '8.30'
$ws.Range('E32').Value = '  +0.92%  '
$ws.Range('E33').Value = '  -2.02%  '
$ws.Range('D34').Value = '3.523.42'
$ws.Range('E34').Value = '  +0.28%  '
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('E36').Value = '  +0.13%  '
$ws.Range('D37').Value = '23.58'
$ws.Range('E37').Value = '  -1.63%  '
$ws.Range('D38').Value = '5.32'
$ws.Range('E38').Value = '  +2.37%  '
$ws.Range('E39').Value = '  -2.74%  '
$ws.Range('D40').Value = '6.87'
$ws.Range('E40').Value = '  -1.49%  '
$ws.Range('D41').Value = '164.23'
$ws.Range('E41').Value = '  -4.29%  '
$ws.Range('D42').Value = '0.0784'
$ws.Range('E42').Value = '  -2.97%  '
$ws.Range('D43').Value = '0.808'
$ws.Range('E43').Value = '  -0.80%  '
$ws.Range('D44').Value = '26.07'
$ws.Range('E44').Value = '  -1.30%  '
$ws.Range('E45').Value = '  +0.03%  '
$ws.Range('E46').Value = '  -0.71%  '
$ws.Range('D47').Value = '41.82'
$ws.Range('E47').Value = '  -1.17%  '
$ws.Range('D48').Value = '4.39'
$ws.Range('E48').Value = '  -0.44%  '
$ws.Range('D49').Value = '1.63'
$ws.Range('E49').Value = '  -1.09%  '
$ws.Range('D50').Value = '2.479.92'
$ws.Range('E50').Value = '  -0.09%  '
$ws.Range('B51').Value = 'SuiNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D51').Value = '0.917'
$ws.Range('E51').Value = '  +1.56%  '
